$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Update the G5:L5 row values (codes changed for this row)
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 1

# Update view state: zoom to 70% and move the selection to E15
$excel.ActiveWindow.Zoom = 70
$ws.Range("E15").Select()
